# Update the "Saldo_guide" workbook:
#  - bump the reference date in column G (Dt. Referencia) from 45488 (2024-07-15)
#    to 45489 (2024-07-16) for every data row
#  - refresh a handful of "Vl. Projetado" / "Vl. Total" pairs (columns E and H)
#    whose totals were recalculated for the new extraction
#  - rename the sheet to reflect the new extraction timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the reference date for every data row (2..275); all of them
# move forward by one day in this extraction.
for ($r = 2; $r -le 275; $r++) {
    $ws.Cells.Item($r, 7).Value = 45489
}

# A few rows also got their projected / total values recalculated.
$updates = @{
    43  = 1526.27
    52  = 2680.09
    55  = 2126.29
    102 = 726.16
    103 = 1904.13
    172 = 849.76
}

foreach ($r in $updates.Keys) {
    $val = $updates[$r]
    $ws.Cells.Item($r, 5).Value = $val
    $ws.Cells.Item($r, 8).Value = $val
}

# New extraction run -> new sheet name.
$ws.Name = "IClientBalance-20240716-093745-"
